$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.713.97'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.82%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.923.43'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.42%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '354.84'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.25'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.61%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.567'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.71%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.630'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.24'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0894'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.10%  '

$ws.Range("E12").Value = '  +0.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.73'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.93'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.386.77'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.927.21'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.979'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.803.27'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.57'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.44%  '

$ws.Range("E20").Value = '  -2.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.11'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0983'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.93'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.84'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.81'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.184'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +12.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.10'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.66%  '

$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.17%  '

$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.51'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +17.71%  '

$ws.Range("E30").Value = '  +13.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.60'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '38.13'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.52%  '

$ws.Range("E33").Value = '  +4.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '52.30'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.59%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0441'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.14%  '

$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("B37").Value = 'Toncoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.89'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -15.90%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.24'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.44'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.60%  '

$ws.Range("E40").Value = '  -0.96%  '

$ws.Range("E41").Value = '  +2.88%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.119'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.03'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.13%  '

$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.86'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.52%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.16'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.46%  '

$ws.Range("E46").Value = '  +2.10%  '

$ws.Range("E47").Value = '  -3.36%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.136.55'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.74%  '

$ws.Range("E49").Value = '  -6.79%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0332'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.54%  '

$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.10'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.87%  '

